# Adds a new "tarifa" (fare) column (E) to every yearly worksheet (2015-2021),
# filling in the header + per-month values, then moves the active
# sheet/selection from sheet "2015" (D15) to sheet "2021" (E15) to match the
# author's final interactive state.

$wb = $excel.ActiveWorkbook

# Per-sheet "tarifa" values for data rows 2..13 (months Jan..Dec).
$values = @{
    1 = @(3.4, 3.4, 3.4, 3.4, 3.4, 3.4, 3.4, 3.4, 3.4, 3.4, 3.4, 3.4)
    2 = @(3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.8)
    3 = @(3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.8, 3.6, 3.6, 3.4, 3.4)
    4 = @(3.4, 3.6, 3.6, 3.6, 3.6, 3.95, 3.95, 3.95, 3.95, 3.95, 3.95, 3.95)
    5 = @(4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05)
    6 = @(4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05)
    7 = @(4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05, 4.05)
}

for ($i = 1; $i -le 7; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $ws.Range("E1").Value = "tarifa"

    $rowValues = $values[$i]
    for ($r = 2; $r -le 13; $r++) {
        $ws.Range("E$r").Value = $rowValues[$r - 2]
    }
}

# Match the author's final selection on every sheet (D15 -> E15), then make
# "2021" (7th tab) the active sheet/selection, mirroring the tabSelected move
# away from "2015" (1st tab) seen in the diff.
for ($i = 1; $i -le 7; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("E15").Select() | Out-Null
}

$ws7 = $wb.Worksheets.Item(7)
$ws7.Activate()
$ws7.Range("E15").Select() | Out-Null
